$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the bold/bordered
# style already used by the other header cells in row 1 (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Style = $ws.Range("H1").Style

# Data values for columns I and J, rows 2-30.
$values = @{
    2  = @(2, 4)
    3  = @(1, 4)
    4  = @(1, 5)
    5  = @(1, 4)
    6  = @(1, 6)
    7  = @(1, 6)
    8  = @(1, 3)
    9  = @(1, 7)
    10 = @(1, 6)
    11 = @(1, 6)
    12 = @(6, 6)
    13 = @(6, 7)
    14 = @(12, 12)
    15 = @(7, 8)
    16 = @(9, 9)
    17 = @(1, 5)
    18 = @(5, 8)
    19 = @(3, 7)
    20 = @(6, 6)
    21 = @(11, 11)
    22 = @(8, 8)
    23 = @(5, 8)
    24 = @(1, 5)
    25 = @(2, 6)
    26 = @(10, 10)
    27 = @(6, 8)
    28 = @(3, 6)
    29 = @(1, 3)
    30 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
